$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("expected")

# Set A3 on "expected" sheet to "index" (new header label for the index column)
$ws2.Range("A3").Value = "index"

# Clear the "Other" text values in B5 and B8 on "expected" sheet
$ws2.Range("B5").Value = ""
$ws2.Range("B8").Value = ""

# Make "expected" sheet the active sheet/tab, with A4 selected
$ws2.Select()
$ws2.Range("A4").Select()

$wb.Save()
